# Glossar.xlsx — add a new glossary entry "Piken" as row 19, pushing the
# existing "Rahmenbedingungen der Detektion" block (rows 19-29) down to
# rows 20-30.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 19..29 down to 20..30, bottom-up so we don't clobber data we
# still need to read. Range.Copy (cell-to-cell, not EntireRow.Insert) keeps
# values + styles without Excel minting a brand-new cellXfs entry.
for ($r = 29; $r -ge 19; $r--) {
    $src = $ws.Range("A" + $r + ":B" + $r)
    $dst = $ws.Range("A" + ($r + 1) + ":B" + ($r + 1))
    $src.Copy($dst)
}

# New row 19: term/definition pair, styled like the other glossary rows
# (copy the formatting from row 2, then overwrite the text).
$ws.Range("A2:B2").Copy($ws.Range("A19:B19"))
$ws.Range("A19").Value = "Piken"
$ws.Range("B19").Value = "Unmittelbar aufeinanderfolgende Ab-Aufwärtsbewegung zur Aufnahme von Zigarettenstummel"

# Match the saved selection state (B19 active cell).
$ws.Range("B19").Select() | Out-Null
